# error solve ifrs list
# Correct the financial figures for rows 2-6 (company 1-5) and remove the
# stray/incorrect trailing rows 7-9 (company 6-8) data cells, keeping only
# the identifying columns A-C there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1242
$ws.Range("E2").Value = 61
$ws.Range("F2").Value = 61
$ws.Range("G2").Value = 84
$ws.Range("H2").Value = 67
$ws.Range("I2").Value = 67
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1234
$ws.Range("L2").Value = 405
$ws.Range("M2").Value = 829
$ws.Range("N2").Value = 827
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 67
$ws.Range("Q2").Value = 109
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = -115
$ws.Range("T2").Value = 10
$ws.Range("U2").Value = 98
$ws.Range("V2").Value = 170
$ws.Range("W2").Value = 4.9
$ws.Range("X2").Value = 5.39
$ws.Range("Y2").Value = 8.42
$ws.Range("Z2").Value = 5.3
$ws.Range("AA2").Value = 48.82
$ws.Range("AB2").Value = 1144.67
$ws.Range("AC2").Value = 702
$ws.Range("AD2").Value = 12.14
$ws.Range("AE2").Value = 8680
$ws.Range("AF2").Value = 0.98
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 2.34
$ws.Range("AI2").Value = 17.05
$ws.Range("AJ2").Value = 9600000

# Row 3
$ws.Range("D3").Value = 1254
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = 65
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 30
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 1370
$ws.Range("L3").Value = 521
$ws.Range("M3").Value = 848
$ws.Range("N3").Value = 847
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 67
$ws.Range("Q3").Value = 26
$ws.Range("R3").Value = -153
$ws.Range("S3").Value = 126
$ws.Range("T3").Value = 155
$ws.Range("U3").Value = -129
$ws.Range("V3").Value = 307
$ws.Range("W3").Value = 5.2
$ws.Range("X3").Value = 2.29
$ws.Range("Y3").Value = 3.53
$ws.Range("Z3").Value = 2.21
$ws.Range("AA3").Value = 61.45
$ws.Range("AB3").Value = 1173.77
$ws.Range("AC3").Value = 308
$ws.Range("AD3").Value = 31.52
$ws.Range("AE3").Value = 8886
$ws.Range("AF3").Value = 1.09
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 2.06
$ws.Range("AI3").Value = 39.03
$ws.Range("AJ3").Value = 9600000

# Row 4
$ws.Range("D4").Value = 1244
$ws.Range("E4").Value = 42
$ws.Range("F4").Value = 42
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 1315
$ws.Range("L4").Value = 468
$ws.Range("M4").Value = 847
$ws.Range("N4").Value = 846
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 67
$ws.Range("Q4").Value = 82
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = -81
$ws.Range("T4").Value = 7
$ws.Range("U4").Value = 75
$ws.Range("V4").Value = 238
$ws.Range("W4").Value = 3.4
$ws.Range("X4").Value = 0.87
$ws.Range("Y4").Value = 1.36
$ws.Range("Z4").Value = 0.8
$ws.Range("AA4").Value = 55.32
$ws.Range("AB4").Value = 1172.1
$ws.Range("AC4").Value = 120
$ws.Range("AD4").Value = 85.59999999999999
$ws.Range("AE4").Value = 8874
$ws.Range("AF4").Value = 1.16
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 1.94
$ws.Range("AI4").Value = 99.67
$ws.Range("AJ4").Value = 9600000

# Row 5
$ws.Range("D5").Value = 1211
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = 40
$ws.Range("G5").Value = 50
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1313
$ws.Range("L5").Value = 449
$ws.Range("M5").Value = 864
$ws.Range("N5").Value = 864
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 67
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = 11
$ws.Range("S5").Value = -37
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = 15
$ws.Range("V5").Value = 213
$ws.Range("W5").Value = 3.33
$ws.Range("X5").Value = 2.62
$ws.Range("Y5").Value = 3.71
$ws.Range("Z5").Value = 2.42
$ws.Range("AA5").Value = 51.98
$ws.Range("AB5").Value = 1197.28
$ws.Range("AC5").Value = 330
$ws.Range("AD5").Value = 33.64
$ws.Range("AE5").Value = 9053
$ws.Range("AF5").Value = 1.23
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 1.8
$ws.Range("AI5").Value = 36.42
$ws.Range("AJ5").Value = 9600000

# Row 6
$ws.Range("D6").Value = 1215
$ws.Range("E6").Value = 11
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = 31
$ws.Range("I6").Value = 31
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = 1327
$ws.Range("L6").Value = 454
$ws.Range("M6").Value = 873
$ws.Range("N6").Value = 873
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 67
$ws.Range("Q6").Value = 20
$ws.Range("R6").Value = 29
$ws.Range("S6").Value = -48
$ws.Range("T6").Value = 36
$ws.Range("U6").Value = -16
$ws.Range("V6").Value = 177
$ws.Range("W6").Value = 0.91
$ws.Range("X6").Value = 2.51
$ws.Range("Y6").Value = 3.52
$ws.Range("Z6").Value = 2.31
$ws.Range("AA6").Value = 51.98
$ws.Range("AB6").Value = 1210.87
$ws.Range("AC6").Value = 318
$ws.Range("AD6").Value = 32.38
$ws.Range("AE6").Value = 9150
$ws.Range("AF6").Value = 1.13
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 1.94
$ws.Range("AI6").Value = 37.64
$ws.Range("AJ6").Value = 9600000

# Rows 7-9: remove all data cells except A, B, C
$ws.Range("D7:AJ9").ClearContents()
